$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.262.97'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '1.676.38'
$ws.Range('E3').Value = '  +0.58%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '217.40'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  +0.40%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '0.5341'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  +4.78%  '
$ws.Range('E7').Value = '  +0.25%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.2678'
$cell.Style = "Normal"
$ws.Range('E8').Value = '  +1.35%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '0.06471'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  +0.00%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.07523'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('D12').Value = '1.675.26'
$ws.Range('E12').Value = '  +0.60%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '4.512'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  +0.33%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '0.5764'
$cell.Style = "Normal"
$ws.Range('E14').Value = '  -1.10%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '0.000008464'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -0.46%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '64.56'
$cell.Style = "Normal"
$ws.Range('D17').Value = '26.266.51'
$ws.Range('E17').Value = '  +0.75%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '4.911'
$cell.Style = "Normal"
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('E21').Value = '  +0.15%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '6.186'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('E23').Value = '  +0.21%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '144.86'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  -0.18%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '0.1278'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  +6.82%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '7.810'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  +2.88%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '15.76'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  +0.97%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '0.06460'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -2.87%  '
$ws.Range('E29').Value = '  +3.77%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '1.318'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  +0.41%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '3.581'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  +1.71%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '3.585'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  +2.17%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '1.654'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  +1.17%  '
$ws.Range('E34').Value = '  +1.30%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '0.6183'
$cell.Style = "Normal"
$ws.Range('E35').Value = '  +1.73%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '2.404'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  +1.56%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '2.722'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  +0.21%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '6.263'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  +0.92%  '
$ws.Range('D39').Value = '1.115.00'
$ws.Range('E39').Value = '  +3.73%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '0.01622'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  +1.31%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.8738'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  +1.51%  '
$ws.Range('E42').Value = '  +0.59%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '100.36'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = '1.827.67'
$ws.Range('E44').Value = '  +0.69%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '0.00000000111'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  -2.99%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '56.93'
$cell.Style = "Normal"
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '8.179'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  +1.73%  '
$ws.Range('E48').Value = '  -0.66%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '0.05259'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  +0.92%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.4288'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  +0.03%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '6.077'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  +2.14%  '
